# Updates cryptocurrency price/volume data on the "cryptos" sheet to reflect
# the latest scrape (commit: "Updated symbol list on Thu Jan 26 11:43:04 UTC 2023
# with GitHub Actions"). Only the Price (column D) and Volume(1h) (column E)
# text values are refreshed; everything else (coin name, link, date, hour,
# styles) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '304.49' },
    @{ Cell = 'E2'; Value = '1.10%' },
    @{ Cell = 'D3'; Value = '35.71' },
    @{ Cell = 'E3'; Value = '1.63%' },
    @{ Cell = 'D4'; Value = '5.070' },
    @{ Cell = 'E4'; Value = '0.55%' },
    @{ Cell = 'D5'; Value = '0.08048' },
    @{ Cell = 'E5'; Value = '1.02%' },
    @{ Cell = 'D6'; Value = '1.918' },
    @{ Cell = 'E6'; Value = '1.38%' },
    @{ Cell = 'D7'; Value = '4.183' },
    @{ Cell = 'E7'; Value = '3.22%' },
    @{ Cell = 'D8'; Value = '7.731' },
    @{ Cell = 'E8'; Value = '-0.86%' },
    @{ Cell = 'D9'; Value = '0.9301' },
    @{ Cell = 'E9'; Value = '0.76%' },
    @{ Cell = 'D10'; Value = '0.1394' },
    @{ Cell = 'E10'; Value = '11.76%' },
    @{ Cell = 'D11'; Value = '0.1898' },
    @{ Cell = 'E11'; Value = '2.37%' },
    @{ Cell = 'D12'; Value = '0.09224' },
    @{ Cell = 'E12'; Value = '-7.37%' },
    @{ Cell = 'D13'; Value = '0.03631' },
    @{ Cell = 'E13'; Value = '1.36%' },
    @{ Cell = 'D14'; Value = '0.09808' },
    @{ Cell = 'E14'; Value = '-0.41%' },
    @{ Cell = 'D15'; Value = '0.001410' },
    @{ Cell = 'E15'; Value = '0.99%' },
    @{ Cell = 'D16'; Value = '0.005904' },
    @{ Cell = 'E16'; Value = '-0.23%' },
    @{ Cell = 'E17'; Value = '1.50%' },
    @{ Cell = 'D18'; Value = '3.015' },
    @{ Cell = 'E18'; Value = '3.47%' },
    @{ Cell = 'E19'; Value = '2.02%' },
    @{ Cell = 'D20'; Value = '0.1304' },
    @{ Cell = 'E20'; Value = '0.21%' },
    @{ Cell = 'D21'; Value = '4.882' },
    @{ Cell = 'E21'; Value = '-3.22%' },
    @{ Cell = 'E22'; Value = '4.53%' },
    @{ Cell = 'D23'; Value = '0.04440' },
    @{ Cell = 'E23'; Value = '-1.29%' },
    @{ Cell = 'E24'; Value = '0.71%' },
    @{ Cell = 'D25'; Value = '0.004784' },
    @{ Cell = 'E25'; Value = '0.01%' },
    @{ Cell = 'E26'; Value = '24.71%' },
    @{ Cell = 'D27'; Value = '0.0003131' },
    @{ Cell = 'E27'; Value = '4.27%' },
    @{ Cell = 'D39'; Value = '0.01961' },
    @{ Cell = 'E39'; Value = '4.19%' },
    @{ Cell = 'D40'; Value = '0.04907' },
    @{ Cell = 'E40'; Value = '4.19%' },
    @{ Cell = 'D41'; Value = '0.007631' },
    @{ Cell = 'E41'; Value = '1.54%' },
    @{ Cell = 'D42'; Value = '0.009258' },
    @{ Cell = 'E42'; Value = '-9.68%' },
    @{ Cell = 'D43'; Value = '0.1374' },
    @{ Cell = 'E44'; Value = '-0.57%' },
    @{ Cell = 'E45'; Value = '6.20%' },
    @{ Cell = 'D46'; Value = '0.00006381' },
    @{ Cell = 'E46'; Value = '2.02%' },
    @{ Cell = 'D47'; Value = '0.00000000751' },
    @{ Cell = 'E47'; Value = '-0.07%' },
    @{ Cell = 'D48'; Value = '65.22' },
    @{ Cell = 'E48'; Value = '1.15%' },
    @{ Cell = 'D49'; Value = '0.001192' },
    @{ Cell = 'E49'; Value = '-20.00%' },
    @{ Cell = 'D50'; Value = '0.00002102' },
    @{ Cell = 'E50'; Value = '-0.07%' },
    @{ Cell = 'D51'; Value = '0.0002002' },
    @{ Cell = 'E51'; Value = '-0.07%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force the cell to be treated as plain text so values such as "304.49"
    # or "1.10%" are stored literally (matching the source inline-string
    # cells) instead of being auto-converted into a number/percentage by
    # Excel. Restoring the style to "Normal" afterwards keeps the original
    # (unstyled) cell formatting intact.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
